$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.735.09'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '3.166.17'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'613.08"
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("D6").Value = "'146.08"
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.160.89'
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = "'35.69"
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = '3.687.34'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D17").Value = '64.708.81'
$ws.Range("D18").Value = '3.166.42'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = "'6.87"
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = "'479.58"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = "'14.61"
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").Value = "'7.92"
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("D24").Value = "'13.75"
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = "'84.12"
$ws.Range("E25").Value = '  +0.44%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  +2.71%  '
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("E29").Value = '  +4.33%  '
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("E31").Value = '  -5.49%  '
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").Value = "'2.68"
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").Value = "'26.56"
$ws.Range("D35").Value = "'1.12"
$ws.Range("E35").Value = '  +0.87%  '
$ws.Range("D36").Value = '0.0₃0790'
$ws.Range("E36").Value = '  +7.09%  '
$ws.Range("D37").Value = "'6.00"
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("D38").Value = "'53.35"
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("D39").Value = "'3.20"
$ws.Range("E39").Value = '  +2.17%  '
$ws.Range("D40").Value = "'460.96"
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  -2.33%  '
$ws.Range("D43").Value = "'8.33"
$ws.Range("E43").Value = '  -1.05%  '
$ws.Range("D44").Value = '2.861.92'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = '  +5.73%  '
$ws.Range("D48").Value = "'26.62"
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = "'35.68"
$ws.Range("E50").Value = '  +6.70%  '
$ws.Range("E51").Value = '  -0.13%  '

# Reset style index back to the default (General) so the quote-prefix
# formatting flag used above does not linger as a visible cell style.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
